# Auto-generated script: restore G-column 'Recorded By' email order
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range('G2').Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range('G3').Value = 'alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg, gehanadel@med.asu.edu.eg, backup@backdoor.com'
$ws.Range('G4').Value = 'asmaa.reda@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hananragab@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, heba@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range('G9').Value = 'Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg'
$ws.Range('G12').Value = 'dalia.t.abdelmegid@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg'
$ws.Range('G14').Value = 'Alaa.A.Mostafa@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg'
$ws.Range('G15').Value = 'Alaa.A.Mostafa@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg'
$ws.Range('G17').Value = 'Walaa.h.ghanima@med.asu.edu.eg, enas.omran@med.asu.edu.eg, user@user.com'
$ws.Range('G18').Value = 'eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range('G19').Value = 'marinasorial@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G20').Value = 'Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G21').Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range('G22').Value = 'alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg, gehanadel@med.asu.edu.eg, backup@backdoor.com'
$ws.Range('G23').Value = 'asmaa.reda@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hananragab@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, heba@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range('G24').Value = 'AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Range('G25').Value = 'yasmintarek@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg'
$ws.Range('G28').Value = 'dina.adel@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, basma.hamed@med.asu.edu.eg'
$ws.Range('G29').Value = 'amira.m.ibrahim@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range('G31').Value = 'dalia.t.abdelmegid@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg'
$ws.Range('G32').Value = 'Noran.Mahmoud@med.asu.edu.eg, menna-allah.gamil@med.asu.edu.eg'
$ws.Range('G34').Value = 'Alaa.A.Mostafa@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg'
$ws.Range('G36').Value = 'Walaa.h.ghanima@med.asu.edu.eg, enas.omran@med.asu.edu.eg, user@user.com'
$ws.Range('G37').Value = 'eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range('G38').Value = 'Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G39').Value = 'Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G40').Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range('G41').Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, hananragab@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range('G42').Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Range('G43').Value = 'AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Amera.a.saad@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg'
$ws.Range('G47').Value = 'nourhan.osama@med.asu.edu.eg, merna.said@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range('G48').Value = 'Eman.m.abosakaya@med.asu.edu.eg, merna.said@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, amany.raafat@med.asu.edu.eg'
$ws.Range('G50').Value = 'dalia.t.abdelmegid@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg'
$ws.Range('G51').Value = 'Noran.Mahmoud@med.asu.edu.eg, menna-allah.gamil@med.asu.edu.eg'
$ws.Range('G52').Value = 'Alaa.A.Mostafa@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg'
$ws.Range('G54').Value = 'Amr-Saeed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg'
$ws.Range('G56').Value = 'eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range('G57').Value = 'Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G58').Value = 'Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G59').Value = 'asmaa.reda@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range('G60').Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, hananragab@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range('G61').Value = 'asmaa.reda@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg'
$ws.Range('G66').Value = 'dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range('G67').Value = 'amira.m.ibrahim@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range('G75').Value = 'eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range('G76').Value = 'marinasorial@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G77').Value = 'Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G78').Value = 'asmaa.reda@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range('G79').Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, hananragab@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range('G80').Value = 'asmaa.reda@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg'
$ws.Range('G81').Value = 'AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Range('G82').Value = 'yasmintarek@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg'
$ws.Range('G83').Value = 'Safa.hany@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Range('G85').Value = 'dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range('G86').Value = 'Eman.m.abosakaya@med.asu.edu.eg, merna.said@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, amany.raafat@med.asu.edu.eg'
$ws.Range('G88').Value = 'dalia.t.abdelmegid@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg'
$ws.Range('G94').Value = 'eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range('G95').Value = 'marinasorial@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G96').Value = 'Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G97').Value = 'asmaa.reda@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range('G98').Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, hananragab@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range('G99').Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Range('G100').Value = 'AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Amera.a.saad@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg'
$ws.Range('G104').Value = 'nourhan.osama@med.asu.edu.eg, merna.said@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range('G113').Value = 'eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range('G115').Value = 'Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
